# Updates cryptocurrency price/volume data on the worksheet to the
# latest scraped values (GitHub Actions refresh of cryptos.xlsx).
# Column D = Price, Column E = Volume(1h) change.
#
# Price values that look numeric are entered with a leading apostrophe
# (quote-prefix) so Excel keeps them as literal text, matching the
# inconsistent "thousands-dot" text formatting already used in the sheet
# (e.g. "63.711.27"), instead of silently converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.711.27"
$ws.Range("E2").Value = "  -3.23%  "
$ws.Range("D3").Value = "2.611.06"
$ws.Range("E3").Value = "  -1.99%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'573.65"
$ws.Range("E5").Value = "  -4.40%  "
$ws.Range("D6").Value = "'155.96"
$ws.Range("E6").Value = "  -2.49%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -3.25%  "
$ws.Range("D9").Value = "2.609.27"
$ws.Range("E9").Value = "  -1.98%  "
$ws.Range("E10").Value = "  -7.17%  "
$ws.Range("D11").Value = "'5.83"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "'0.382"
$ws.Range("E12").Value = "  -4.65%  "
$ws.Range("D13").Value = "'0.156"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").Value = "'28.12"
$ws.Range("E14").Value = "  -3.70%  "
$ws.Range("D15").Value = "3.078.07"
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("E16").Value = "  -7.99%  "
$ws.Range("D17").Value = "63.593.68"
$ws.Range("D18").Value = "2.611.94"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("D19").Value = "'11.98"
$ws.Range("E19").Value = "  -4.84%  "
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("E21").Value = "  -6.00%  "
$ws.Range("D22").Value = "'342.54"
$ws.Range("E22").Value = "  -3.83%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  -3.69%  "
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("E26").Value = "  -4.69%  "
$ws.Range("D27").Value = "'589.86"
$ws.Range("E27").Value = "  +3.28%  "
$ws.Range("E28").Value = "  -6.61%  "
$ws.Range("E29").Value = "  -3.37%  "
$ws.Range("E30").Value = "  -1.94%  "
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("E32").Value = "  -3.12%  "
$ws.Range("E33").Value = "  -3.61%  "
$ws.Range("E34").Value = "  -6.02%  "
$ws.Range("D35").Value = "'6.60"
$ws.Range("E35").Value = "  -1.97%  "
$ws.Range("D36").Value = "'5.42"
$ws.Range("E36").Value = "  -1.45%  "
$ws.Range("D37").Value = "'0.402"
$ws.Range("E37").Value = "  -5.17%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  -4.45%  "
$ws.Range("D40").Value = "'154.08"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  -5.91%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "'2.51"
$ws.Range("E43").Value = "  +1.39%  "
$ws.Range("E44").Value = "  -3.40%  "
$ws.Range("D45").Value = "'156.80"
$ws.Range("D46").Value = "'23.73"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("D47").Value = "'3.89"
$ws.Range("E47").Value = "  -5.37%  "
$ws.Range("E48").Value = "  -4.93%  "
$ws.Range("D49").Value = "'0.629"
$ws.Range("E49").Value = "  -2.42%  "
$ws.Range("D50").Value = "'0.100"
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("D51").Value = "'0.0248"
$ws.Range("E51").Value = "  -4.48%  "
